# add a test scenequest. support mechanical move back
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new quest entry ("losttree;2") to the QuestDungeon list for the
# first dungeon row (row 4 / Id 18000001, the "losttrees" test scene).
$cell = $ws.Range("L4")
$cell.Value = $cell.Value2 + "|losttree;2"

# Move selection back onto the cell we just edited.
$ws.Range("L4").Select()
